# Switch the deck's applied design back from the custom "Integral" (Red
# Violet) theme to the built-in "Office Theme".
#
# Real PowerPoint would do this by picking "Office Theme" in the Design
# gallery, which rewrites the colour scheme carried by the slide master's
# theme part. We reproduce that colour-scheme rewrite through the
# Theme/ThemeColorScheme object model.

$p = $ppt.ActivePresentation

# --- Re-point every theme colour at the stock Office palette --------------
# Order matches OOXML a:clrScheme: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink. ThemeColor.RGB uses the usual OLE 0xBBGGRR layout, so build each
# value from its normal RRGGBB hex form to avoid mistakes.
function Get-OleRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $themeColors.Colors($i + 1).RGB = Get-OleRgb $officeColors[$i]
}
